# Applies corrected IFRS figures (per-share/per-unit financials) to rows 2-9
# of the company_list sheet, replacing the previously mis-scaled values and
# blanking a handful of cells that have no data for the corrected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 34379
$ws.Range("E2").Value = 3584
$ws.Range("F2").Value = 3584
$ws.Range("G2").Value = 1734
$ws.Range("H2").Value = 1316
$ws.Range("I2").Value = 1296
$ws.Range("J2").Value = 21
$ws.Range("K2").Value = 48047
$ws.Range("L2").Value = 34787
$ws.Range("M2").Value = 13260
$ws.Range("N2").Value = 12090
$ws.Range("O2").Value = 1170
$ws.Range("P2").Value = 7900
$ws.Range("Q2").Value = 4998
$ws.Range("R2").Value = -2926
$ws.Range("S2").Value = -1732
$ws.Range("T2").Value = 2828
$ws.Range("U2").Value = 2170
$ws.Range("V2").Value = 24216
$ws.Range("W2").Value = 10.42
$ws.Range("X2").Value = 3.83
$ws.Range("Y2").Value = 11.64
$ws.Range("Z2").Value = 2.76
$ws.Range("AA2").Value = 262.35
$ws.Range("AB2").Value = 51.12
$ws.Range("AC2").Value = 860
$ws.Range("AD2").Value = 11.24
$ws.Range("AE2").Value = 7652
$ws.Range("AF2").Value = 1.26
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 157993158

# Row 3
$ws.Range("D3").Value = 30404
$ws.Range("E3").Value = 1360
$ws.Range("F3").Value = 1360
$ws.Range("G3").Value = -655
$ws.Range("H3").Value = -675
$ws.Range("I3").Value = -694
$ws.Range("J3").Value = 19
$ws.Range("K3").Value = 52200
$ws.Range("L3").Value = 39592
$ws.Range("M3").Value = 12608
$ws.Range("N3").Value = 11429
$ws.Range("O3").Value = 1179
$ws.Range("P3").Value = 7900
$ws.Range("Q3").Value = 2250
$ws.Range("R3").Value = -5817
$ws.Range("S3").Value = 3907
$ws.Range("T3").Value = 6758
$ws.Range("U3").Value = -4507
$ws.Range("V3").Value = 27805
$ws.Range("W3").Value = 4.47
$ws.Range("X3").Value = -2.22
$ws.Range("Y3").Value = -5.9
$ws.Range("Z3").Value = -1.35
$ws.Range("AA3").Value = 314.02
$ws.Range("AB3").Value = 41.96
$ws.Range("AC3").Value = -439
$ws.Range("AD3").Value = -15.33
$ws.Range("AE3").Value = 7234
$ws.Range("AF3").Value = 0.93
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 157993158

# Row 4
$ws.Range("D4").Value = 29472
$ws.Range("E4").Value = 1201
$ws.Range("F4").Value = 1201
$ws.Range("G4").Value = -378
$ws.Range("H4").Value = -379
$ws.Range("I4").Value = -360
$ws.Range("J4").Value = -18
$ws.Range("K4").Value = 51217
$ws.Range("L4").Value = 39076
$ws.Range("M4").Value = 12141
$ws.Range("N4").Value = 11022
$ws.Range("O4").Value = 1119
$ws.Range("P4").Value = 7900
$ws.Range("Q4").Value = 1883
$ws.Range("R4").Value = -2649
$ws.Range("S4").Value = 56
$ws.Range("T4").Value = 3655
$ws.Range("U4").Value = -1771
$ws.Range("V4").Value = 27331
$ws.Range("W4").Value = 4.07
$ws.Range("X4").Value = -1.29
$ws.Range("Y4").Value = -3.21
$ws.Range("Z4").Value = -0.73
$ws.Range("AA4").Value = 321.85
$ws.Range("AB4").Value = 38.06
$ws.Range("AC4").Value = -228
$ws.Range("AD4").Value = -37.17
$ws.Range("AE4").Value = 6976
$ws.Range("AF4").Value = 1.22
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 157993158

# Row 5
$ws.Range("D5").Value = 28764
$ws.Range("E5").Value = -1572
$ws.Range("F5").Value = -1572
$ws.Range("G5").Value = -1318
$ws.Range("H5").Value = -1118
$ws.Range("I5").Value = -1248
$ws.Range("J5").Value = 130
$ws.Range("K5").Value = 45009
$ws.Range("L5").Value = 35095
$ws.Range("M5").Value = 9914
$ws.Range("N5").Value = 9545
$ws.Range("O5").Value = 369
$ws.Range("P5").Value = 7900
$ws.Range("Q5").Value = 713
$ws.Range("R5").Value = -837
$ws.Range("S5").Value = 49
$ws.Range("T5").Value = 2899
$ws.Range("U5").Value = -2185
$ws.Range("V5").Value = 24042
$ws.Range("W5").Value = -5.46
$ws.Range("X5").Value = -3.89
$ws.Range("Y5").Value = -12.14
$ws.Range("Z5").Value = -2.32
$ws.Range("AA5").Value = 353.99
$ws.Range("AB5").Value = 22.7
$ws.Range("AC5").Value = -790
$ws.Range("AD5").Value = -5.6
$ws.Range("AE5").Value = 6042
$ws.Range("AF5").Value = 0.73
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 157993158

# Row 6
$ws.Range("D6").Value = 25587
$ws.Range("E6").Value = -789
$ws.Range("F6").Value = -789
$ws.Range("G6").Value = -1817
$ws.Range("H6").Value = -1827
$ws.Range("I6").Value = -1742
$ws.Range("K6").Value = 42147
$ws.Range("L6").Value = 28018
$ws.Range("M6").Value = 14130
$ws.Range("N6").Value = 14130
$ws.Range("P6").Value = 14363
$ws.Range("Q6").Value = 568
$ws.Range("R6").Value = -1672
$ws.Range("S6").Value = 2288
$ws.Range("T6").Value = 1588
$ws.Range("U6").Value = -1020
$ws.Range("V6").Value = 19362
$ws.Range("W6").Value = -3.08
$ws.Range("X6").Value = -7.14
$ws.Range("Y6").Value = -14.72
$ws.Range("Z6").Value = -4.19
$ws.Range("AA6").Value = 198.29
$ws.Range("AB6").Value = -0.92
$ws.Range("AC6").Value = -787
$ws.Range("AD6").Value = -6.75
$ws.Range("AE6").Value = 4919
$ws.Range("AF6").Value = 1.08
$ws.Range("AJ6").Value = 287260287
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").ClearContents()

# Row 7
$ws.Range("D7").Value = 23815
$ws.Range("E7").Value = 250
$ws.Range("G7").Value = -780
$ws.Range("H7").Value = -765
$ws.Range("I7").Value = -765
$ws.Range("K7").Value = 39890
$ws.Range("L7").Value = 26815
$ws.Range("M7").Value = 13080
$ws.Range("N7").Value = 13520
$ws.Range("P7").Value = 14360
$ws.Range("Q7").Value = 675
$ws.Range("R7").Value = -1075
$ws.Range("S7").Value = -1775
$ws.Range("T7").Value = 1360
$ws.Range("W7").Value = 1.05
$ws.Range("X7").Value = -3.21
$ws.Range("Y7").Value = -5.53
$ws.Range("Z7").Value = -1.87
$ws.Range("AA7").Value = 205.01
$ws.Range("AC7").Value = -266
$ws.Range("AD7").Value = -14.34
$ws.Range("AE7").Value = 4707
$ws.Range("AF7").Value = 0.8100000000000001
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("U7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 24220
$ws.Range("E8").Value = 420
$ws.Range("G8").Value = -425
$ws.Range("H8").Value = -410
$ws.Range("I8").Value = -410
$ws.Range("K8").Value = 39660
$ws.Range("L8").Value = 26865
$ws.Range("M8").Value = 12800
$ws.Range("N8").Value = 13255
$ws.Range("P8").Value = 14360
$ws.Range("Q8").Value = 1890
$ws.Range("R8").Value = -1215
$ws.Range("S8").Value = 110
$ws.Range("T8").Value = 1390
$ws.Range("W8").Value = 1.73
$ws.Range("X8").Value = -1.69
$ws.Range("Y8").Value = -3.06
$ws.Range("Z8").Value = -1.03
$ws.Range("AA8").Value = 209.88
$ws.Range("AC8").Value = -143
$ws.Range("AD8").Value = -26.76
$ws.Range("AE8").Value = 4614
$ws.Range("AF8").Value = 0.83
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("U8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 25190
$ws.Range("E9").Value = 590
$ws.Range("G9").Value = -255
$ws.Range("H9").Value = -245
$ws.Range("I9").Value = -245
$ws.Range("K9").Value = 40055
$ws.Range("L9").Value = 27375
$ws.Range("M9").Value = 12680
$ws.Range("N9").Value = 13150
$ws.Range("P9").Value = 14360
$ws.Range("Q9").Value = 1810
$ws.Range("R9").Value = -1390
$ws.Range("S9").Value = 345
$ws.Range("T9").Value = 1415
$ws.Range("W9").Value = 2.34
$ws.Range("X9").Value = -0.97
$ws.Range("Y9").Value = -1.86
$ws.Range("Z9").Value = -0.62
$ws.Range("AA9").Value = 215.89
$ws.Range("AC9").Value = -85
$ws.Range("AD9").Value = -44.79
$ws.Range("AE9").Value = 4578
$ws.Range("AF9").Value = 0.83
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("U9").ClearContents()
$ws.Range("AI9").ClearContents()
